$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

# ---------------------------------------------------------------------------
# 1) "- [Visual Studio](...)" paragraph: merge split runs, drop proofErr tags,
#    and add the es-CO language tag to the paragraph mark + every run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(4)
$xml = @"
<w:p $wNs>
  <w:pPr>
    <w:rPr>
      <w:lang w:val="es-CO"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-CO"/>
    </w:rPr>
    <w:t>- [Visual Studio](</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
      <w:lang w:val="es-CO"/>
    </w:rPr>
    <w:t>https://visualstudio.microsoft.com/</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="es-CO"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 2) "- [.NET 8.0](...)" paragraph: merge split runs, drop proofErr tags.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(5)
$xml = @"
<w:p $wNs>
  <w:r>
    <w:t>- [.NET 8.0](</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>https://dotnet.microsoft.com/en-us/download/dotnet/8.0</w:t>
  </w:r>
  <w:r>
    <w:t>)</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) "- [MSSQL Server](...)" paragraph: merge split runs, drop proofErr tags,
#    replace the fldChar HYPERLINK field with a real w:hyperlink element.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(6)
$xml = @"
<w:p $wNs>
  <w:r>
    <w:t>- [MSSQL Server](</w:t>
  </w:r>
  <w:hyperlink r:id="rId5" w:history="1">
    <w:r>
      <w:rPr>
        <w:rStyle w:val="Hyperlink"/>
      </w:rPr>
      <w:t>https://www.microsoft.com/en-us/sql-server/sql-server-downloads</w:t>
    </w:r>
  </w:hyperlink>
  <w:r>
    <w:t>)</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) "- [SQL Server Management Studio](...)" paragraph: merge split runs,
#    drop proofErr tags.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(7)
$xml = @"
<w:p $wNs>
  <w:r>
    <w:t>- [SQL Server Management Studio](</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
    </w:rPr>
    <w:t>https://learn.microsoft.com/en-us/sql/ssms/download-sql-server-management-studio-ssms?view=sql-server-ver16#download-ssms</w:t>
  </w:r>
  <w:r>
    <w:t>)</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 5) "Create a new MSSQL server using and make sure..." paragraph: merge
#    split runs, drop proofErr tags.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(10)
$xml = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Create a new </w:t>
  </w:r>
  <w:r>
    <w:t>MSSQL server</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> using and make sure that it is running</w:t>
  </w:r>
  <w:r>
    <w:t>. You can check on the SQL configuration manager</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------------------
# 6) "Paste the connection string on the Appsettings.json file" paragraph:
#    merge split runs, drop proofErr/spellErr tags, drop trailing colon, and
#    insert a brand-new "NOTE: ..." paragraph right after it (picking up the
#    lastRenderedPageBreak that used to sit on the following image run).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(18)
$xml = @"
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Paste the connection string on the Appsettings.json file</w:t>
  </w:r>
</w:p>
<w:p $wNs>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>NOTE: Do not forget to change &#8220;master&#8221; by your database name</w:t>
  </w:r>
  <w:r>
    <w:t>. Also, add TrustServerCertificate=True to your connection string.</w:t>
  </w:r>
</w:p>
"@
$p.Range.InsertXML($xml)

# Remove the now-redundant lastRenderedPageBreak from the following image run
# (it moved to the new NOTE paragraph above).
$p = $d.Paragraphs(20)
Write-Output $p.Range.Text
$p.Range.Find.Execute("", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
